$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three sensor ids that were removed from the working set.
# Delete bottom-up so earlier row numbers stay valid as later ones are removed.
$ws.Rows("29").Delete()   # 77067
$ws.Rows("27").Delete()   # 128349
$ws.Rows("7").Delete()    # 124737

# Column B only ever held a stray, empty styled cell (B2) - remove it entirely.
$ws.Columns("B").Delete()

# The remaining ids no longer carry the leftover "applyFont" style.
$ws.Range("A2:A27").ClearFormats()

# Shrink the conditional-formatting range by the same 3 rows that were removed.
$cf = $ws.Range("A31:A46").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A28:A43"))

# Shrink the hidden _FilterDatabase name to match the new data extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$A`$43"
    }
}

# Move the active selection to B12, matching where editing left off.
$ws.Range("B12").Select()
